$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corrections (Appenzeller-Herzog (2019) - van Dis (2020))
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

$ws.Range("H3").Value = 0.8261582615826158
$ws.Range("I3").Value = 0.01497534921939195
$ws.Range("J3").Value = 0.9
$ws.Range("K3").Value = 78.40000000000001

$ws.Range("Q3").Value = 20
$ws.Range("R3").Value = 53
$ws.Range("S3").Value = 65
$ws.Range("T3").Value = 78
$ws.Range("U3").Value = 108

$ws.Range("V3").Value = 4848
$ws.Range("W3").Value = 4815
$ws.Range("X3").Value = 4803
$ws.Range("Y3").Value = 4790
$ws.Range("Z3").Value = 4760

$ws.Range("AF3").Value = 0.995892
$ws.Range("AG3").Value = 0.989113
$ws.Range("AH3").Value = 0.9866470000000001
$ws.Range("AI3").Value = 0.983977
$ws.Range("AJ3").Value = 0.977814

$wb.Save()
